# Apply the edit described by the diff:
# - A new record row is inserted at row 142 (all existing rows from 142
#   downward shift down by one, ending at row 161).
# - The new row 142 contains a new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 142; this pushes the former rows
# 142..160 down to 143..161 (values, formats & styles move with them).
$ws.Rows.Item(142).Insert()

# Populate the newly inserted row 142 with the new data record.
$ws.Range("A142").Value2 = 4
$ws.Range("B142").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C142").Value2 = "Los Lagos"
$ws.Range("D142").Value2 = 44946
$ws.Range("E142").Value2 = 10
$ws.Range("F142").Value2 = 100112052
$ws.Range("G142").Value2 = "Albahaca"
$ws.Range("H142").Value2 = "Sin especificar"
$ws.Range("I142").Value2 = "Primera"
$ws.Range("J142").Value2 = 90
$ws.Range("K142").Value2 = 6000
$ws.Range("L142").Value2 = 6000
$ws.Range("M142").Value2 = 6000
$ws.Range("N142").Value2 = "$/docena de matas"
$ws.Range("O142").Value2 = "Región Metropolitana"
$ws.Range("P142").Value2 = 1000
$ws.Range("Q142").Value2 = 6
$ws.Range("R142").Value2 = "Hortaliza"
